$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new data rows (102-146) following the same repeating pattern as
# the existing rows: regcntr_id cycles 10002..10010, device_id increments by 1,
# lang_code="eng", is_active=TRUE, cr_by="superadmin", cr_dtimes="now()".
$regIds = 10002,10003,10004,10005,10006,10007,10008,10009,10010
$startDevice = 3000121
$startRow = 102
$endRow = 146

for ($row = $startRow; $row -le $endRow; $row++) {
    $offset = $row - $startRow
    $regId = $regIds[$offset % $regIds.Length]
    $deviceId = $startDevice + $offset

    $ws.Cells.Item($row, 1).Value = $regId
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Reflect the new selection/scroll state that Excel recorded after entering
# the additional rows (active cell A102, selection spanning the new rows).
[void]$ws.Range("A102:F146").Select()

# Page was set up for printing (portrait orientation) as part of this edit.
$ws.PageSetup.Orientation = 1
